# Physical system is ready for sending data.
# Extends Sheet1 with 15 new columns (T:AH, i.e. columns 20-34) of
# sensor/measurement data for each of the existing rows 1-6, mirroring the
# existing B:S pattern already on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 20  # column T

# --- Row 1 (bold/bordered header-style numeric row) ---------------------
# B1:S1 already carry the bordered/centered style; stamp the same format
# onto the new cells before writing their values.
$ws.Range("B1").Copy($ws.Range("T1:AH1"))
$row1 = @(27,28,29,30,34,35,36,38,39,40,41,43,42,44,45)
for ($i = 0; $i -lt $row1.Count; $i++) {
    $ws.Cells.Item(1, $firstCol + $i).Value = $row1[$i]
}

# --- Row 2 (carID) --------------------------------------------------------
$row2 = @(27,28,30,30,34,35,36,38,39,40,43,43,43,44,45)
for ($i = 0; $i -lt $row2.Count; $i++) {
    $ws.Cells.Item(2, $firstCol + $i).Value = $row2[$i]
}

# --- Row 3 (speed2) --------------------------------------------------------
$row3 = @(50,80.62,42.72,75,40,69.64,80.62,79.06,47.43,73.81999999999999,63.64,55.9,55.9,87.45999999999999,69.45999999999999)
for ($i = 0; $i -lt $row3.Count; $i++) {
    $ws.Cells.Item(3, $firstCol + $i).Value = $row3[$i]
}

# --- Row 4 (asma) -----------------------------------------------------------
# Stored as text in the source data (values like "76.70" / "149.90" need
# the trailing zero preserved, which a numeric cell would drop), so force
# text via a leading apostrophe on each cell. That quote-prefix also stamps
# a "quotePrefix" style onto the cell, so immediately strip it back down to
# the default (unstyled) format via a formats-only paste from a bare cell,
# same as the rest of row 4.
$row4 = @("42.86","130.35","22.06","114.29","14.29","98.98","130.35","125.88","35.53","110.93","81.83","59.72","59.72","149.90","98.46")
for ($i = 0; $i -lt $row4.Count; $i++) {
    $ws.Cells.Item(4, $firstCol + $i).Formula = "'" + $row4[$i]
}
$ws.Range("A1").Copy()
$ws.Range("T4:AH4").PasteSpecial(-4122)

# --- Row 5 (ceza_tutar) ------------------------------------------------------
$row5 = @(3136,6440,1508.5,6440,1508.5,6440,6440,6440,3136,6440,6440,6440,6440,6440,6440)
for ($i = 0; $i -lt $row5.Count; $i++) {
    $ws.Cells.Item(5, $firstCol + $i).Value = $row5[$i]
}

# --- Row 6 (hesaplanan_asma) -------------------------------------------------
$row6 = @(30,50,10,50,10,50,50,50,30,50,50,50,50,50,50)
for ($i = 0; $i -lt $row6.Count; $i++) {
    $ws.Cells.Item(6, $firstCol + $i).Value = $row6[$i]
}
